# Applies the gh-pages data refresh (output generated at 456a3b4)
# Updates "想去人数" (F) and "最低票价" (G) columns across the
# 展览 / 演出 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 455
$ws.Range("G5").Value = 58
$ws.Range("G6").Value = 35
$ws.Range("F7").Value = 578
$ws.Range("F8").Value = 80
$ws.Range("G8").Value = "不可售"
$ws.Range("F9").Value = 6896
$ws.Range("G9").Value = 80
$ws.Range("F16").Value = 16364
$ws.Range("F17").Value = 4
$ws.Range("F23").Value = 11453
$ws.Range("F25").Value = 1091
$ws.Range("F26").Value = 4509
$ws.Range("F27").Value = 366

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = 180

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 455
$ws.Range("G5").Value = 58
$ws.Range("G6").Value = 35
$ws.Range("F7").Value = 578
$ws.Range("F9").Value = 80
$ws.Range("G9").Value = "不可售"
$ws.Range("F10").Value = 6896
$ws.Range("G10").Value = 80
$ws.Range("F18").Value = 16364
$ws.Range("F19").Value = 4
$ws.Range("G25").Value = 180
$ws.Range("F27").Value = 11453
$ws.Range("F29").Value = 1091
$ws.Range("F30").Value = 4509
$ws.Range("F31").Value = 366
